# ------------------------------------------------------------------
# Applies the commit: adds a new "Player Info" sheet as the first
# sheet in the workbook, and reworks the MATCH_CARD_LINK column in
# the "ODI Batting" and "ODI Bowling" sheets into a MATCH_CODE column
# that stores just the numeric HowStat match code instead of the
# full scorecard URL.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Player Info" worksheet in front of "ODI Batting"
# ------------------------------------------------------------------
$battingSheetBefore = $wb.Worksheets.Item("ODI Batting")

$playerInfo = $wb.Worksheets.Add($battingSheetBefore)
$playerInfo.Name = "Player Info"

# NOTE: worksheet references captured before the insertion above can
# become stale (the engine seems to resolve held sheet references by
# position rather than identity), so re-fetch the sheets we still
# need to edit by name, now that the sheet collection is stable.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRng = $playerInfo.Range("A1:D1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# Data row - keep the player id as text, matching the rest of the workbook
$dataRng = $playerInfo.Range("A2:D2")
$dataRng.NumberFormat = "@"
$playerInfo.Range("A2").Value = "5926"
$playerInfo.Range("B2").Value = "Cameron Green"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"
$dataRng.Style = "Normal"

# ------------------------------------------------------------------
# 2. Rework MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" (column D)
#    and "ODI Bowling" (column B): keep only the numeric match code
#    instead of the full scorecard URL.
# ------------------------------------------------------------------

# --- ODI Batting: MATCH_CARD_LINK lives in column D ---
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"
$lastRow = $battingSheet.UsedRange.Rows.Count
$dataRange = $battingSheet.Range("D2:D$lastRow")
$dataRange.NumberFormat = "@"
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Value2
    $pieces = $link -split "MatchCode="
    $cell.Value = $pieces[1]
}
$dataRange.Style = "Normal"

# --- ODI Bowling: MATCH_CARD_LINK lives in column B ---
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"
$lastRow2 = $bowlingSheet.UsedRange.Rows.Count
$dataRange2 = $bowlingSheet.Range("B2:B$lastRow2")
$dataRange2.NumberFormat = "@"
for ($r = 2; $r -le $lastRow2; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Value2
    $pieces = $link -split "MatchCode="
    $cell.Value = $pieces[1]
}
$dataRange2.Style = "Normal"

Write-Output "done"
